$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "outliers" column (L) moves to before "Account type" (J),
# pushing "Account type" and "Account hierarchy" one column to the right
# (J->K, K->L). Equivalent to Excel's "Cut column, then Insert Cut Cells".
$ws.Columns("L").Cut()
$ws.Columns("J").Insert()
